# Actualización desde MV -datos-
# Appends the new daily-rate rows (15-09-2021 .. 30-09-2021) to Sheet1,
# right after the existing last row (178, dated 14-09-2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Serie (date label), 3 meses, 6 meses, 1 año
$rows = @(
    @("15-09-2021", -0.02, 0.33, 0.39),
    @("16-09-2021", -0.14, 0.29, 0.38),
    @("20-09-2021", -0.18, 0.30, 0.39),
    @("21-09-2021", -0.18, 0.39, 0.49),
    @("22-09-2021", -0.27, 0.24, 0.42),
    @("23-09-2021", -0.47, 0.13, 0.37),
    @("24-09-2021", -0.64, 0.03, 0.43),
    @("27-09-2021", -0.69, -0.07, 0.43),
    @("28-09-2021", -0.61, -0.09, 0.21),
    @("29-09-2021", -0.90, -0.06, 0.20),
    @("30-09-2021", -0.99, -0.31, 0.13)
)

$startRow = 179
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
